$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update the four input values that drive the downstream formulas.
$ws.Range("A40").Value = 3
$ws.Range("A41").Value = 2
$ws.Range("A42").Value = 2
$ws.Range("A43").Value = 3

# Move the active selection to C49, matching the saved view state.
$ws.Activate()
$ws.Range("C49").Select()
